# Edit: "Upcoming Meeting.pptx" supervisor-meeting update
#
# 1. Delete the slide with sldId 264 (slide index 7) - a minimal
#    "Research Questions" slide that only had the Title, a rounded
#    rectangle and a "Research Question 3" label left on it.
# 2. On the slide that was pushed up into position 7 (sldId 265, the
#    big "Research Questions" slide), remove the nine leftover
#    Research-Question-1/2/3 shapes (ids 8,10,13,14,15,16,17,18,19),
#    keeping only the title, the "When thinking about robust..."
#    textbox and its rounded-rectangle background.
# 3. Nudge the flowchart picture on slide 4 (sldId 262) down slightly
#    (y offset 1733752 EMU -> 1815945 EMU).

$p = $ppt.ActivePresentation

# --- 1. Remove the slide that only has the RQ3 leftovers (sldId 264) ---
$p.Slides.Item(7).Delete()

# --- 2. Clean up the shapes on the slide that is now at position 7 ---
$s = $p.Slides.Item(7)
$idsToDelete = @(8, 10, 13, 14, 15, 16, 17, 18, 19)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($idsToDelete -contains $shp.Id) {
        $shp.Delete()
    }
}

# --- 3. Move the flowchart picture on slide 4 down a little ---
$picSlide = $p.Slides.Item(4)
for ($i = 1; $i -le $picSlide.Shapes.Count; $i++) {
    $shp = $picSlide.Shapes.Item($i)
    if ($shp.Name -eq "Picture 3") {
        $shp.Top = 142.98785
    }
}
